$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F23").Value = "So Akira knew Misaki was here..."
$ws.Range("G24").Value = "When I look at him, he has a genuinely surprised look on his face."
$ws.Range("F28").Value = "Now that I think about it, is Nanase-kun too?"
$ws.Range("G29").Value = "We're on our way home now, but..."
$ws.Range("F30").Value = "Ah, that's right."
$ws.Range("F36").Value = "But today is different."
$ws.Range("F37").Value = "Ah-... I failed..."
$ws.Range("F38").Value = "Can't you just buy the magazine?"
$ws.Range("F39").Value = "But it's expensive. And if it's a special that features an author I don't like, wouldn't it be a waste?"
$ws.Range("G41").Value = "That's just how it is."
$ws.Range("G43").Value = "It's your fault for forgetting. Let's do it tomorrow."
$ws.Range("G48").Value = "I also mention that he forgot that new books were arriving at the library (while trying to make it sound interesting)."
$ws.Range("F50").Value = "Wasn't there a special on Gardner this month? (Gardner is a lawyer and author)"
$ws.Range("H50").Value = "Deliberately add context"
$ws.Range("F52").Value = "Wasn't he the guy who created the fictional defense attorney, Perry Mason?"
$ws.Range("H52").Value = "Deliberately add context"
$ws.Range("G54").Value = "I see. Fufu. Good thing in the midst of misfortune...?"
$ws.Range("H54").Value = "くすっ is a giggle."
$ws.Range("G58").Value = "Akira just laughs simply. He really is a small happiness man."
$ws.Range("H58").Value = ""
$ws.Range("F61").Value = "...Well, never mind."
$ws.Range("F62").Value = "That's why, you see, Akira, you better hurry or you'll be late."
$ws.Range("F66").Value = "I didn't, I didn't. Well, I can hang out and then go home though."
$ws.Range("G68").Value = "Fufu... Fujii-kun..."
$ws.Range("H68").Value = "くすっ is a giggle"
$ws.Range("G73").Value = "Go with Misaki-san to the library."

$ws.Range("G12").Select()
